$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the latest reporting/publish date column (H) with the new database entry
$ws.Range("H9").Value = "1402-01-28 (8)"

# Updated read_price algorithm now yields real figures instead of placeholder dashes
$ws.Range("H25").Value = 0
$ws.Range("H26").Value = 7580
